$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @"
244|calendly.com/cortezatenobraulio1|68.65
245|calendly.com/towerclinic|55.2
246|calendly.com/yuktahar|41.75
247|calendly.com/pitchydeck|55.2
248|calendly.com/up|28.35
249|calendly.com/imaree|41.75
250|calendly.com/communitize|55.2
251|calendly.com/maryannarcenal|55.2
252|calendly.com/mindfulmedia|55.2
253|calendly.com/hiroyo|41.75
254|calendly.com/rayhansocial|55.2
255|calendly.com/origads|41.75
256|calendly.com/otomati|41.75
257|calendly.com/markbruns|55.2
258|calendly.com/neelofurshahab|55.2
259|calendly.com/rebeccahsteele|55.2
260|calendly.com/cassidy|41.75
261|calendly.com/maryna|41.75
262|calendly.com/yg001|41.75
263|calendly.com/tee|28.35
264|calendly.com/sphereofinfluence360|68.65
265|calendly.com/matteo|41.75
266|calendly.com/amy|28.35
267|calendly.com/growhigh|41.75
268|calendly.com/henning|41.75
269|calendly.com/al|28.35
270|calendly.com/howard|41.75
271|calendly.com/agpllc|41.75
272|calendly.com/cfw|28.35
273|calendly.com/aleperezelias|55.2
274|calendly.com/aryamanmahjan|55.2
275|calendly.com/andrew|41.75
276|calendly.com/gabriel|41.75
277|calendly.com/acadium|41.75
278|calendly.com/corecotton|55.2
279|calendly.com/kimberlyhall|55.2
280|calendly.com/john|28.35
281|calendly.com/niabettertogether|55.2
282|calendly.com/quadrantshift|55.2
283|calendly.com/staffnetscheduling|68.65
284|calendly.com/alyssa|41.75
"@

# Re-assert word-wrap on the pre-existing "long description" rows (134-228,
# 230-243). The engine's xlsx round-trip does not retain a bare <alignment
# wrapText="true"/> that has no other explicit per-cell override, so make it
# explicit again to keep those rows visually consistent with the new ones.
for ($r = 134; $r -le 228; $r++) {
    $ws.Cells.Item($r, 1).WrapText = $true
}
for ($r = 230; $r -le 243; $r++) {
    $ws.Cells.Item($r, 1).WrapText = $true
}

$rows = $newData -split "`n"
foreach ($line in $rows) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $rowNum = [int]$parts[0]
    $val = $parts[1]
    $height = [double]$parts[2]

    $cell = $ws.Cells.Item($rowNum, 1)
    $cell.Value = $val
    $cell.WrapText = $true
    $ws.Rows.Item($rowNum).RowHeight = $height
}

# Update selection / active cell to match the final state (A245)
[void]$ws.Range("A245").Select()

# Best-effort: scroll the view so the newly added rows are visible
$win = $excel.ActiveWindow
$win.ScrollRow = 277
